$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter the new "Phân công" (assignment) names in the same order they were
# first introduced so the shared-string table gets the same new entries in
# the same order (Tú, Nhi -> Huy, Kiều -> Lê, Kiều).
$ws.Range("F8").Value = "Tú, Nhi"
$ws.Range("F6").Value = "Huy, Kiều"
$ws.Range("F34").Value = "Lê, Kiều"

# Fill in the matching "Hoàn thành" (completion) percentages.
$ws.Range("E6").Value = 0.7
$ws.Range("E8").Value = 1
$ws.Range("E34").Value = 0.9

# Update the current selection to reflect where the user ended up.
$ws.Activate() | Out-Null
$ws.Range("F35").Select() | Out-Null

Write-Output "done"
